# Added OHC user as required for rota meeting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing staff-list values (column B, rows 6-12) before
# shifting them down to make room for the new "MahaDeva AM" entry.
$b6  = $ws.Range("B6").Value2
$b7  = $ws.Range("B7").Value2
$b8  = $ws.Range("B8").Value2
$b9  = $ws.Range("B9").Value2
$b10 = $ws.Range("B10").Value2
$b11 = $ws.Range("B11").Value2
$b12 = $ws.Range("B12").Value2

# Rename the AM/PM split for the existing "MahaDeva" header entry
$ws.Range("A3").Value2 = "MahaDeva PM"

# New OHC attendee cell, previously empty
$ws.Range("C3").Value2 = "MahaDeva - OHC"

# Insert "MahaDeva AM" into the staff list (column B), pushing the
# existing names down by one row starting at row 6
$ws.Range("B6").Value2  = "MahaDeva AM"
$ws.Range("B7").Value2  = $b6
$ws.Range("B8").Value2  = $b7
$ws.Range("B9").Value2  = $b8
$ws.Range("B10").Value2 = $b9
$ws.Range("B11").Value2 = $b10
$ws.Range("B12").Value2 = $b11
$ws.Range("B13").Value2 = $b12

# Update the summary note in the footer row
$ws.Range("F35").Value2 = "MahaDeva Mahi Ben "
